$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6549
$ws.Range("I62").Value = 3145.75
$ws.Range("J62").Value = 8250.625
$ws.Range("K62").Value = 3145.75
$ws.Range("L62").Value = 8250.625
$ws.Range("M62").Value = -2521.75
$ws.Range("N62").Value = -9498.625
$ws.Range("H65").Value = 6549
$ws.Range("I65").Value = 3145.75
$ws.Range("J65").Value = 8250.625
$ws.Range("K65").Value = 15728.75
$ws.Range("L65").Value = 41253.125
$ws.Range("M65").Value = -12608.75
$ws.Range("N65").Value = -47493.125
$ws.Range("H96").Value = 1471
$ws.Range("J96").Value = 375
$ws.Range("L96").Value = 1125
$ws.Range("N96").Value = -3871
$ws.Range("H116").Value = 7702.7334
$ws.Range("I116").Value = 6803.273
$ws.Range("K116").Value = 6803.273
$ws.Range("M116").Value = -3361.273
$ws.Range("H132").Value = 1896.9642
$ws.Range("I132").Value = 1795.2307
$ws.Range("J132").Value = 3219.5
$ws.Range("K132").Value = 5385.6921
$ws.Range("L132").Value = 9658.5
$ws.Range("M132").Value = -2855.6921
$ws.Range("N132").Value = -14718.5
$ws.Range("H135").Value = 5078
$ws.Range("I135").Value = 4941.143
$ws.Range("J135").Value = 6036
$ws.Range("K135").Value = 44470.287
$ws.Range("L135").Value = 54324
$ws.Range("M135").Value = -41935.287
$ws.Range("N135").Value = -59394
$ws.Range("H138").Value = 3598.4707
$ws.Range("I138").Value = 1640.2858
$ws.Range("J138").Value = 3774.205
$ws.Range("K138").Value = 4920.857400000001
$ws.Range("L138").Value = 11322.615
$ws.Range("M138").Value = 219.1425999999992
$ws.Range("N138").Value = -21602.615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3960.7058
$ws.Range("I2").Value = 1151.3
$ws.Range("K2").Value = 1151.3
$ws.Range("M2").Value = -1038.3
$ws.Range("H31").Value = 27615.5
$ws.Range("I31").Value = 15231.333
$ws.Range("J31").Value = 39999.668
$ws.Range("K31").Value = 15231.333
$ws.Range("L31").Value = 39999.668
$ws.Range("M31").Value = -14937.333
$ws.Range("N31").Value = -40587.668
$ws.Range("H32").Value = 8930.777
$ws.Range("I32").Value = 7555.706
$ws.Range("J32").Value = 13181
$ws.Range("K32").Value = 7555.706
$ws.Range("L32").Value = 13181
$ws.Range("M32").Value = -7268.706
$ws.Range("N32").Value = -13755
$ws.Range("H43").Value = 35666.332
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 35666.332
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 35666.332
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -36292.332
$ws.Range("H110").Value = 2081.5898
$ws.Range("I110").Value = 1614.1818
$ws.Range("K110").Value = 1614.1818
$ws.Range("M110").Value = 430.8181999999999
$ws.Range("H116").Value = 3960.7058
$ws.Range("I116").Value = 1151.3
$ws.Range("K116").Value = 1151.3
$ws.Range("M116").Value = 1142.7
$ws.Range("H122").Value = 2805.8647
$ws.Range("I122").Value = 2375.92
$ws.Range("K122").Value = 7127.76
$ws.Range("M122").Value = -4677.76
$ws.Range("H132").Value = 2306.3555
$ws.Range("I132").Value = 1790.4359
$ws.Range("K132").Value = 5371.307699999999
$ws.Range("M132").Value = -2841.307699999999
$ws.Range("H139").Value = 66412
$ws.Range("J139").Value = 66412
$ws.Range("L139").Value = 66412
$ws.Range("N139").Value = -76692

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3960.7058
$ws.Range("I3").Value = 1151.3
$ws.Range("K3").Value = 1151.3
$ws.Range("M3").Value = -1037.3
$ws.Range("H105").Value = 37519.555
$ws.Range("I105").Value = 50750
$ws.Range("K105").Value = 50750
$ws.Range("M105").Value = -49003
$ws.Range("H107").Value = 1846.3529
$ws.Range("I107").Value = 1899.25
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1899.25
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 20.75
$ws.Range("N107").Value = -4840
$ws.Range("H134").Value = 2756.9429
$ws.Range("I134").Value = 2199.3103
$ws.Range("K134").Value = 6597.9309
$ws.Range("M134").Value = -4062.9309

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2885.3076
$ws.Range("I16").Value = 2100.9
$ws.Range("K16").Value = 2100.9
$ws.Range("M16").Value = -1813.9
$ws.Range("H31").Value = 28255.121
$ws.Range("I31").Value = 1717.0385
$ws.Range("K31").Value = 1717.0385
$ws.Range("M31").Value = -1422.0385
$ws.Range("H34").Value = 28255.121
$ws.Range("I34").Value = 1717.0385
$ws.Range("K34").Value = 1717.0385
$ws.Range("M34").Value = -1515.0385
$ws.Range("H41").Value = 1059
$ws.Range("I41").Value = 1059
$ws.Range("K41").Value = 1059
$ws.Range("M41").Value = -631
$ws.Range("H95").Value = 32063.75
$ws.Range("J95").Value = 32063.75
$ws.Range("L95").Value = 32063.75
$ws.Range("N95").Value = -37555.75
$ws.Range("H113").Value = 2885.3076
$ws.Range("I113").Value = 2100.9
$ws.Range("K113").Value = 2100.9
$ws.Range("M113").Value = 69.09999999999991
$ws.Range("H134").Value = 35716940
$ws.Range("I134").Value = 47621170
$ws.Range("K134").Value = 142863510
$ws.Range("M134").Value = -142860975
$ws.Range("H141").Value = 227009.47
$ws.Range("J141").Value = 237824.44
$ws.Range("L141").Value = 237824.44
$ws.Range("N141").Value = -248184.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 136.3
$ws.Range("J12").Value = 191.57143
$ws.Range("L12").Value = 574.71429
$ws.Range("N12").Value = -920.71429
$ws.Range("H134").Value = 4748.385
$ws.Range("I134").Value = 4748.385
$ws.Range("K134").Value = 14245.155
$ws.Range("M134").Value = -9175.155000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5386.273
$ws.Range("J43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("N43").Value = -25302
$ws.Range("H80").Value = 839082.7
$ws.Range("I80").Value = 1668663.4
$ws.Range("K80").Value = 1668663.4
$ws.Range("M80").Value = -1667665.4
$ws.Range("H83").Value = 839082.7
$ws.Range("I83").Value = 1668663.4
$ws.Range("K83").Value = 8343317
$ws.Range("M83").Value = -8338325
$ws.Range("H93").Value = 38267.332
$ws.Range("J93").Value = 38267.332
$ws.Range("L93").Value = 38267.332
$ws.Range("N93").Value = -42011.332
$ws.Range("H122").Value = 3092.6667
$ws.Range("I122").Value = 2666.6333
$ws.Range("J122").Value = 3802.7222
$ws.Range("K122").Value = 7999.8999
$ws.Range("L122").Value = 11408.1666
$ws.Range("M122").Value = -5549.8999
$ws.Range("N122").Value = -16308.1666
$ws.Range("H132").Value = 4811451.5
$ws.Range("I132").Value = 6412651.5
$ws.Range("K132").Value = 19237954.5
$ws.Range("M132").Value = -19235424.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 21272.727
$ws.Range("J2").Value = 24666.666
$ws.Range("L2").Value = 24666.666
$ws.Range("N2").Value = -24890.666
$ws.Range("H46").Value = 2874.9167
$ws.Range("I46").Value = 1812.375
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1812.375
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -1624.375
$ws.Range("N46").Value = -5376
$ws.Range("H55").Value = 1472730.4
$ws.Range("I55").Value = 2632672.5
$ws.Range("J55").Value = 3470.2
$ws.Range("K55").Value = 2632672.5
$ws.Range("L55").Value = 3470.2
$ws.Range("M55").Value = -2632499.5
$ws.Range("N55").Value = -3816.2
$ws.Range("H68").Value = 8976.388999999999
$ws.Range("I68").Value = 6796.7144
$ws.Range("K68").Value = 6796.7144
$ws.Range("M68").Value = -6047.7144
$ws.Range("H71").Value = 8976.388999999999
$ws.Range("I71").Value = 6796.7144
$ws.Range("K71").Value = 33983.572
$ws.Range("M71").Value = -30239.572
$ws.Range("H82").Value = 1118.32
$ws.Range("I82").Value = 874.7778
$ws.Range("J82").Value = 1255.3125
$ws.Range("K82").Value = 874.7778
$ws.Range("L82").Value = 1255.3125
$ws.Range("M82").Value = -513.7778
$ws.Range("N82").Value = -1977.3125
$ws.Range("H85").Value = 1118.32
$ws.Range("I85").Value = 874.7778
$ws.Range("J85").Value = 1255.3125
$ws.Range("K85").Value = 874.7778
$ws.Range("L85").Value = 1255.3125
$ws.Range("M85").Value = 373.2222
$ws.Range("N85").Value = -3751.3125
$ws.Range("H93").Value = 1503.5416
$ws.Range("I93").Value = 1819.0588
$ws.Range("J93").Value = 737.2857
$ws.Range("K93").Value = 1819.0588
$ws.Range("L93").Value = 737.2857
$ws.Range("M93").Value = -571.0588
$ws.Range("N93").Value = -3233.2857
$ws.Range("H100").Value = 4945
$ws.Range("I100").Value = 4083.3635
$ws.Range("K100").Value = 4083.3635
$ws.Range("M100").Value = -3542.3635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 61710.332
$ws.Range("J109").Value = 61710.332
$ws.Range("L109").Value = 61710.332
$ws.Range("N109").Value = -64484.332
$ws.Range("H132").Value = 13688.25
$ws.Range("I132").Value = 1568.9395
$ws.Range("J132").Value = 147000.67
$ws.Range("K132").Value = 4706.818499999999
$ws.Range("L132").Value = 441002.01
$ws.Range("M132").Value = -2176.818499999999
$ws.Range("N132").Value = -446062.01
$ws.Range("H136").Value = 3007.7144
$ws.Range("I136").Value = 1491.826
$ws.Range("J136").Value = 9980.799999999999
$ws.Range("K136").Value = 4475.478
$ws.Range("L136").Value = 29942.4
$ws.Range("M136").Value = -1925.478
$ws.Range("N136").Value = -35042.39999999999
